# Voice-to-Text requirements sheet update
# - Renames requirement IDs from RWS000xx to RWVT000xx and rewrites their
#   descriptions to reflect the new voice_to_text() function requirements.
# - Drops the old RWS00009..RWS00017 rows (now unused / out of scope).
# - Widens column D, wraps text + grows row 6 for the long multi-line
#   requirement description.
# - Updates the active selection, matching the author's last cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the 8 remaining requirement rows (IDs first, then descriptions,
#     so the shared-string table lays out in the same order as the target) ---
$ws.Range("D2").Value2 = "RWVT00001"
$ws.Range("D3").Value2 = "RWVT00002"
$ws.Range("D4").Value2 = "RWVT00003"
$ws.Range("D5").Value2 = "RWVT00004"
$ws.Range("D6").Value2 = "RWVT00005"
$ws.Range("D7").Value2 = "RWVT00006"
$ws.Range("D8").Value2 = "RWVT00007"
$ws.Range("D9").Value2 = "RWVT00008"

$ws.Range("E2").Value2 = "function should be voice_to_text()"
$ws.Range("E3").Value2 = "Input argument must be JSON"
$ws.Range("E4").Value2 = "the JSON object must contain audio file data, which will be processed by the respective voice-to-text service."
$ws.Range("E5").Value2 = "Return argument must for JSON"
$ws.Range("E6").Value2 = "Return object for JSON have comprised with following`na. status: SUCCESS/ERROR`nb. error: AUDIO_ERROR/API_ERROR/SERVER_BUSY/UNDEFINED`nc. response: the transcribed text from the voice-to-text service"
$ws.Range("E7").Value2 = "the function should establish a connection with a Google voice-to-text service and process the given audio data."
$ws.Range("E8").Value2 = "once the connection is established, the function should send the audio data to the service and retrieve the text response."
$ws.Range("E9").Value2 = "The function must return the JSON object with the transcribed text and relevant status/error details, as specified in RWS00005"

# --- Remove the now-obsolete requirement rows 10-18 (D & E columns) ---
$ws.Range("D10:E18").ClearContents()

# --- Column D is now wider to fit the longer RWVT ids ---
$ws.Columns.Item(4).ColumnWidth = 16.14

# --- Wrap the long multi-line requirement text and grow its row ---
$ws.Range("E6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 86.4

# --- Match the author's final selection ---
$ws.Range("E10").Select() | Out-Null
